# Applies the "Solicitudes de 09_01 y correcciones de 09_02" edit:
#  - Updates several image-request rows (IMG02, IMG05, IMG07) on the
#    "Solicitud gráfica" sheet with corrected / new descriptions.
#  - Fills in previously-empty "Tipo"/"Formato" columns for a couple of rows.
#  - Adds two brand-new image requests (IMG09, IMG10) in rows 18-19.
#  - Updates the workbook view (zoom level + selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solicitud gráfica")

# --- Row 11 (IMG02): was "Hombre tomándose la cabeza..." -> now the
#     "Persona en traslúcido..." photo (note the leading space Aula
#     Planeta typed before "Persona" in the long description path).
$ws.Range("B11").Value = "3° ESO/Biología y geología/La coordinación: sistema nervioso y endocrino/Evaluación/ Persona en traslúcido de color fucsia y en donde se ve el cerebro y los nervios"
$ws.Range("J11").Value = "Persona en traslúcido de color fucsia y en donde se ve el cerebro y los nervios"

# --- Row 13 (IMG04): observation text corrected / shortened.
$ws.Range("K13").Value = "Esta imagen hace parte del recurso ""Evaluación"" del guion ubicado en 3° ESO/Biología y geología/La coordinación: docrino"

# --- Row 14 (IMG05): old "Persona en traslúcido..." request replaced by a
#     brand-new Shutterstock photo reference (numeric id), new description,
#     and the observation is cleared.
$ws.Range("B14").Value = 79216378
$ws.Range("J14").Value = "Fotografía de neurona y células gliales en fondo azul"
$ws.Range("K14").Value = ""

# --- Row 16 (IMG07): Tipo / Formato were blank, now set like the other rows.
$ws.Range("D16").Value = "Fotografía"
$ws.Range("E16").Value = "Horizontal"

# --- Row 17 (IMG08): Tipo / Formato filled in, and an observation added.
$ws.Range("D17").Value = "Fotografía"
$ws.Range("E17").Value = "Horizontal"
$ws.Range("K17").Value = "Esta imagen hace parte del recurso ""Evaluación"" del guion ubicado en 3° ESO/Biología y geología/La coordinación: docrino"

# --- Row 18 (new IMG09 request).
$ws.Range("B18").Value = 311314832
$ws.Range("D18").Value = "Fotografía"
$ws.Range("E18").Value = "Horizontal"
$ws.Range("J18").Value = "Neurona motora y músculo"

# --- Row 19 (new IMG10 request). A19 is typed over with the literal
#     "IMG10" instead of being left as the auto-numbering formula.
$ws.Range("A19").Value = "IMG10"
$ws.Range("B19").Value = 285141302
$ws.Range("D19").Value = "Fotografía"
$ws.Range("E19").Value = "Horizontal"
$ws.Range("J19").Value = "Beirponas café en fondo gris"

# --- View state: zoomed out, and the last touched cell is F19.
$ws.Range("A9").Select()
$excel.ActiveWindow.Zoom = 65
$ws.Range("F19").Select()
